$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 16:52"

# Refresh country statistics and re-sort rows by total cases (column B) descending
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 337925
$ws.Cells.Item(4, 3).Value = 1252
$ws.Cells.Item(4, 4).Value = 18002
$ws.Cells.Item(4, 5).Value = 310259
$ws.Cells.Item(4, 6).Value = 8702
$ws.Cells.Item(4, 7).Value = 48
$ws.Cells.Item(4, 8).Value = 9664

$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 100315
$ws.Cells.Item(7, 3).Value = 192
$ws.Cells.Item(7, 4).Value = 28700
$ws.Cells.Item(7, 5).Value = 70018
$ws.Cells.Item(7, 6).Value = 3936
$ws.Cells.Item(7, 7).Value = 13
$ws.Cells.Item(7, 8).Value = 1597

$ws.Cells.Item(15, 1).Value = "Paises Bajos"
$ws.Cells.Item(15, 2).Value = 18803
$ws.Cells.Item(15, 3).Value = 952
$ws.Cells.Item(15, 4).Value = 250
$ws.Cells.Item(15, 5).Value = 16686
$ws.Cells.Item(15, 6).Value = 1409
$ws.Cells.Item(15, 7).Value = 101
$ws.Cells.Item(15, 8).Value = 1867

$ws.Cells.Item(21, 1).Value = "Israel"
$ws.Cells.Item(21, 2).Value = 8611
$ws.Cells.Item(21, 3).Value = 181
$ws.Cells.Item(21, 4).Value = 585
$ws.Cells.Item(21, 5).Value = 7970
$ws.Cells.Item(21, 6).Value = 141
$ws.Cells.Item(21, 7).Value = 7
$ws.Cells.Item(21, 8).Value = 56

$ws.Cells.Item(48, 1).Value = "Republica Dominicana"
$ws.Cells.Item(48, 2).Value = 1828
$ws.Cells.Item(48, 3).Value = 83
$ws.Cells.Item(48, 4).Value = 33
$ws.Cells.Item(48, 5).Value = 1709
$ws.Cells.Item(48, 6).Value = 147
$ws.Cells.Item(48, 7).Value = 4
$ws.Cells.Item(48, 8).Value = 86

$ws.Cells.Item(64, 1).Value = "Moldavia"
$ws.Cells.Item(64, 2).Value = 965
$ws.Cells.Item(64, 3).Value = 101
$ws.Cells.Item(64, 4).Value = 37
$ws.Cells.Item(64, 5).Value = 909
$ws.Cells.Item(64, 6).Value = 80
$ws.Cells.Item(64, 7).Value = 4
$ws.Cells.Item(64, 8).Value = 19

$ws.Cells.Item(65, 1).Value = "Irak"
$ws.Cells.Item(65, 2).Value = 961
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 279
$ws.Cells.Item(65, 5).Value = 621
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 61

$ws.Cells.Item(66, 1).Value = "Hong Kong"
$ws.Cells.Item(66, 2).Value = 915
$ws.Cells.Item(66, 3).Value = 24
$ws.Cells.Item(66, 4).Value = 216
$ws.Cells.Item(66, 5).Value = 695
$ws.Cells.Item(66, 6).Value = 12
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 4

$ws.Cells.Item(80, 1).Value = "Bulgaria"
$ws.Cells.Item(80, 2).Value = 549
$ws.Cells.Item(80, 3).Value = 18
$ws.Cells.Item(80, 4).Value = 39
$ws.Cells.Item(80, 5).Value = 488
$ws.Cells.Item(80, 6).Value = 26
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 22

$ws.Cells.Item(81, 1).Value = "Letonia"
$ws.Cells.Item(81, 2).Value = 542
$ws.Cells.Item(81, 3).Value = 9
$ws.Cells.Item(81, 4).Value = 16
$ws.Cells.Item(81, 5).Value = 525
$ws.Cells.Item(81, 6).Value = 5
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 1

$ws.Cells.Item(83, 1).Value = "Eslovaquia"
$ws.Cells.Item(83, 2).Value = 534
$ws.Cells.Item(83, 3).Value = 49
$ws.Cells.Item(83, 4).Value = 8
$ws.Cells.Item(83, 5).Value = 524
$ws.Cells.Item(83, 6).Value = 3
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = 2

$ws.Cells.Item(109, 1).Value = "Georgia"
$ws.Cells.Item(109, 2).Value = 188
$ws.Cells.Item(109, 3).Value = 14
$ws.Cells.Item(109, 4).Value = 39
$ws.Cells.Item(109, 5).Value = 147
$ws.Cells.Item(109, 6).Value = 6
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 2

$ws.Cells.Item(179, 1).Value = "Seychelles"
$ws.Cells.Item(179, 2).Value = 11
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 5).Value = 11
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = "Groenlandia"
$ws.Cells.Item(180, 2).Value = 11
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 3
$ws.Cells.Item(180, 5).Value = 8
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "Curazao"
$ws.Cells.Item(181, 2).Value = 11
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 5
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 1

$ws.Cells.Item(182, 1).Value = "Suazilandia"
$ws.Cells.Item(182, 2).Value = 10
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 10
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

$ws.Cells.Item(183, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(183, 2).Value = 10
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 10
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = "Surinam"
$ws.Cells.Item(184, 2).Value = 10
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(185, 1).Value = "Mozambique"
$ws.Cells.Item(185, 2).Value = 10
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 1

$ws.Cells.Item(186, 1).Value = "Republica del Chad"
$ws.Cells.Item(186, 2).Value = 9
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 0

